$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (A385, date-format style)
# down through the new date cells A386:A464 so the appended rows keep the
# same numFmt/border/font/alignment as the rest of column A.
$ws.Range("A385").Copy($ws.Range("A386:A464"))

# New daily figures, continuing the series through 2021-12-08 ("aggiornamento fino a 8/12").
$data = @(
    @(386, 44460, 10, 114, 60.31331178278742),
    @(387, 44461, 3, 101, 53.43547798299588),
    @(388, 44462, 18, 114, 60.31331178278742),
    @(389, 44463, 22, 106, 56.08079867522339),
    @(390, 44464, 4, 94, 49.73202901387736),
    @(391, 44465, 26, 95, 50.26109315232286),
    @(392, 44466, 9, 92, 48.67390073698634),
    @(393, 44467, 4, 86, 45.49951590631333),
    @(394, 44468, 2, 85, 44.97045176786782),
    @(395, 44469, 20, 87, 46.02858004475883),
    @(396, 44470, 22, 87, 46.02858004475883),
    @(397, 44471, 13, 96, 50.79015729076836),
    @(398, 44472, 10, 80, 42.3251310756403),
    @(399, 44473, 21, 92, 48.67390073698634),
    @(400, 44474, 0, 88, 46.55764418320432),
    @(401, 44475, 11, 97, 51.31922142921386),
    @(402, 44476, 21, 98, 51.84828556765937),
    @(403, 44477, 14, 90, 47.61577246009534),
    @(404, 44478, 8, 85, 44.97045176786782),
    @(405, 44479, 13, 88, 46.55764418320432),
    @(406, 44480, 4, 71, 37.56355382963077),
    @(407, 44481, 1, 72, 38.09261796807627),
    @(408, 44482, 5, 66, 34.91823313740325),
    @(409, 44483, 9, 54, 28.5694634760572),
    @(410, 44484, 5, 45, 23.80788623004767),
    @(411, 44485, 12, 49, 25.92414278382969),
    @(412, 44486, 3, 39, 20.63350139937465),
    @(413, 44487, 7, 42, 22.22069381471116),
    @(414, 44488, 7, 48, 25.39507864538418),
    @(415, 44489, 9, 52, 27.51133519916619),
    @(416, 44490, 1, 44, 23.27882209160216),
    @(417, 44491, 8, 47, 24.86601450693868),
    @(418, 44492, 7, 42, 22.22069381471116),
    @(419, 44493, 7, 46, 24.33695036849317),
    @(420, 44494, 6, 45, 23.80788623004767),
    @(421, 44495, 4, 42, 22.22069381471116),
    @(422, 44496, 4, 37, 19.57537312248364),
    @(423, 44497, 5, 41, 21.69162967626566),
    @(424, 44498, 5, 38, 20.10443726092914),
    @(425, 44499, 10, 41, 21.69162967626566),
    @(426, 44500, 9, 43, 22.74975795315666),
    @(427, 44501, 13, 50, 26.45320692227519),
    @(428, 44502, 3, 49, 25.92414278382969),
    @(429, 44503, 0, 45, 23.80788623004767),
    @(430, 44504, 3, 43, 22.74975795315666),
    @(431, 44505, 19, 57, 30.15665589139371),
    @(432, 44506, 25, 72, 38.09261796807627),
    @(433, 44507, 20, 83, 43.91232349097681),
    @(434, 44508, 18, 88, 46.55764418320432),
    @(435, 44509, 6, 91, 48.14483659854084),
    @(436, 44510, 6, 97, 51.31922142921386),
    @(437, 44511, 26, 120, 63.48769661346046),
    @(438, 44512, 15, 116, 61.37144005967844),
    @(439, 44513, 8, 99, 52.37734970610487),
    @(440, 44514, 12, 91, 48.14483659854084),
    @(441, 44515, 16, 89, 47.08670832164984),
    @(442, 44516, 136, 219, 115.8650463195653),
    @(443, 44517, 2, 215, 113.7487897657833),
    @(444, 44518, 18, 207, 109.5162766582193),
    @(445, 44519, 31, 223, 117.9813028733473),
    @(446, 44520, 13, 228, 120.6266235655748),
    @(447, 44521, 11, 227, 120.0975594271294),
    @(448, 44522, 44, 255, 134.9113553036035),
    @(449, 44523, 7, 126, 66.66208144413348),
    @(450, 44524, 131, 255, 134.9113553036035),
    @(451, 44525, 36, 273, 144.4345097956225),
    @(452, 44526, 15, 257, 135.9694835804945),
    @(453, 44527, 34, 278, 147.07983048785),
    @(454, 44528, 44, 311, 164.5389470565516),
    @(455, 44529, 29, 296, 156.6029849798691),
    @(456, 44530, 20, 309, 163.4808187796607),
    @(457, 44531, 11, 189, 99.99312216620021),
    @(458, 44532, 30, 183, 96.81873733552719),
    @(459, 44533, 55, 223, 117.9813028733473),
    @(460, 44534, 42, 231, 122.2138159809114),
    @(461, 44535, 28, 215, 113.7487897657833),
    @(462, 44536, 55, 241, 127.5044573653664),
    @(463, 44537, 24, 245, 129.6207139191484),
    @(464, 44538, 7, 241, 127.5044573653664)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
